# The document ends with a placeholder paragraph containing a single
# ellipsis ("…") followed by a trailing empty paragraph. This edit:
#   1. Splits the ellipsis paragraph into a new blank paragraph followed
#      by a paragraph containing the real sign-off text "Aina Jung.".
#   2. Marks "Aina" with the proofErr spell-check bookends Word adds for
#      words it doesn't recognise, and splits the text into the two runs
#      ("Aina" / " Jung.") that straddle that markup.
#   3. Removes the paragraph break that is now redundant (the document
#      keeps exactly one trailing blank paragraph, same as before).

$d = $word.ActiveDocument

# --- Step 1: turn "…" into a blank line followed by "Aina Jung." -----
# Using Find/Replace with the "^p" paragraph-mark code (rather than
# Range.InsertParagraphBefore) keeps the newly created blank paragraph
# completely empty in the OOXML (no stray empty run left behind).
$d.Content.Find.Execute("…", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "^pAina Jung.", 2)

# --- Step 2: locate the new "Aina Jung." paragraph and inject the ----
# proofErr / run markup via InsertXML (this replaces only the matched
# range's contents, leaving the paragraph's own properties untouched).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Aina Jung.") {
        $target = $p
        break
    }
}

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Aina</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Jung.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$contentRange = $d.Range($target.Range.Start, $target.Range.End - 1)
$contentRange.InsertXML($xml)

# --- Step 3: drop the now-duplicated trailing blank paragraph --------
# The split in step 1 left an extra blank paragraph at the end of the
# document (on top of the one that was already there). Re-resolve the
# "Aina Jung." paragraph (positions shifted after InsertXML) and delete
# the paragraph mark right after it so it merges back down to a single
# trailing blank paragraph, matching the original paragraph count.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Aina Jung.") {
        $target = $p
        break
    }
}
$markRange = $d.Range($target.Range.End - 1, $target.Range.End)
$markRange.Delete()
